$d = $word.ActiveDocument

# --- Paragraph 2: split "Prolog-Interpreter (32bit) ist installiert"
# into two runs: "Prolog-Interpreter" and " ist installiert" ---

# Locate the end of "Prolog-Interpreter" to create a run-boundary there
# (adding a point-bookmark inside a run forces the engine to split it).
$findRng = $d.Paragraphs.Item(2).Range
$findRng.Find.Execute("Prolog-Interpreter", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($findRng.End, $findRng.End)
$d.Bookmarks.Add("_TMP_SPLIT", $splitPoint)

# Remove the "(32bit) " portion, leaving the leading space before "ist".
$para2Again = $d.Paragraphs.Item(2).Range
$para2Again.Find.Execute("(32bit) ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Remove the temporary split-helper bookmark (keeps the run split it created).
$d.Bookmarks.Item("_TMP_SPLIT").Delete()

# --- Move the _GoBack bookmark from the end of paragraph 2 to the start
# of paragraph 1 (collapsed / zero-length), matching the target diff. ---

$d.Bookmarks.Item("_GoBack").Delete()

# A collapsed bookmark placed with Start = End = 0 (the very beginning of the
# document) is mishandled when created directly, so insert a throwaway
# placeholder character at position 0, anchor the collapsed bookmark right
# after it (a non-zero start), then delete the placeholder. The bookmark
# correctly shifts back to position 0 once the character in front of it is
# removed.
$insertPoint = $d.Range(0, 0)
$insertPoint.InsertBefore("X")

$goBackRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)

$d.Range(0, 1).Delete()
